# Update district level data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 6
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 6

# Row 3
$ws.Range("F3").Value = 3

# Row 4
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 3

# Row 5
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 5

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 3

# Update selection to match final state (F6 selected)
$ws.Range("F6").Select()
